$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 240 (existing rows 240:257 shift down to 241:258)
$ws.Rows.Item(240).Insert()

# Populate the new row 240 with the weekly pineapple price record
$ws.Cells.Item(240, 1).Value = 11
$ws.Cells.Item(240, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(240, 3).Value = "Bíobío"
$ws.Cells.Item(240, 4).Value = 45013
$ws.Cells.Item(240, 5).Value = 8
$ws.Cells.Item(240, 6).Value = "Fruta"
$ws.Cells.Item(240, 7).Value = 100108
$ws.Cells.Item(240, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(240, 9).Value = 100108005
$ws.Cells.Item(240, 10).Value = "Piña"
$ws.Cells.Item(240, 11).Value = "Caramelo"
$ws.Cells.Item(240, 12).Value = "Segunda"
$ws.Cells.Item(240, 13).Value = 200
$ws.Cells.Item(240, 14).Value = 20000
$ws.Cells.Item(240, 15).Value = 21000
$ws.Cells.Item(240, 16).Value = 20500
$ws.Cells.Item(240, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(240, 18).Value = "Ecuador"
$ws.Cells.Item(240, 19).Value = 1464
$ws.Cells.Item(240, 20).Value = 14
